$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 38067  # F2: 38053 -> 38067
$ws1.Cells.Item(3, 6).Value = 50  # F3: 49 -> 50
$ws1.Cells.Item(7, 6).Value = 382  # F7: 381 -> 382
$ws1.Cells.Item(9, 6).Value = 872  # F9: 871 -> 872
$ws1.Cells.Item(11, 6).Value = 776  # F11: 771 -> 776
$ws1.Cells.Item(12, 6).Value = 605  # F12: 601 -> 605
$ws1.Cells.Item(13, 6).Value = 97  # F13: 94 -> 97
$ws1.Cells.Item(14, 6).Value = 39  # F14: 40 -> 39
$ws1.Cells.Item(15, 6).Value = 49  # F15: 48 -> 49
$ws1.Cells.Item(16, 6).Value = 704  # F16: 701 -> 704
$ws1.Cells.Item(17, 6).Value = 198  # F17: 197 -> 198
$ws1.Cells.Item(18, 6).Value = 501  # F18: 500 -> 501
$ws1.Cells.Item(20, 6).Value = 1203  # F20: 1201 -> 1203
$ws1.Cells.Item(22, 6).Value = 902  # F22: 898 -> 902
$ws1.Cells.Item(23, 6).Value = 2630  # F23: 2624 -> 2630
$ws1.Cells.Item(24, 6).Value = 1120  # F24: 1111 -> 1120
$ws1.Cells.Item(25, 6).Value = 592  # F25: 590 -> 592
$ws1.Cells.Item(26, 6).Value = 134  # F26: 133 -> 134
$ws1.Cells.Item(29, 6).Value = 867  # F29: 863 -> 867
$ws1.Cells.Item(31, 6).Value = 1206  # F31: 1202 -> 1206

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 471  # F3: 470 -> 471
$ws2.Cells.Item(4, 6).Value = 340  # F4: 339 -> 340
$ws2.Cells.Item(12, 6).Value = 16  # F12: 14 -> 16

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 685  # F2: 684 -> 685

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 685  # F2: 684 -> 685
$ws4.Cells.Item(3, 6).Value = 38067  # F3: 38053 -> 38067
$ws4.Cells.Item(4, 6).Value = 50  # F4: 49 -> 50
$ws4.Cells.Item(9, 6).Value = 382  # F9: 381 -> 382
$ws4.Cells.Item(11, 6).Value = 471  # F11: 470 -> 471
$ws4.Cells.Item(12, 6).Value = 340  # F12: 339 -> 340
$ws4.Cells.Item(15, 6).Value = 872  # F15: 871 -> 872
$ws4.Cells.Item(17, 6).Value = 776  # F17: 771 -> 776
$ws4.Cells.Item(18, 6).Value = 605  # F18: 601 -> 605
$ws4.Cells.Item(19, 6).Value = 98  # F19: 94 -> 98
$ws4.Cells.Item(21, 6).Value = 39  # F21: 40 -> 39
$ws4.Cells.Item(25, 6).Value = 49  # F25: 48 -> 49
$ws4.Cells.Item(27, 6).Value = 704  # F27: 701 -> 704
$ws4.Cells.Item(28, 6).Value = 198  # F28: 197 -> 198
$ws4.Cells.Item(29, 6).Value = 501  # F29: 500 -> 501
$ws4.Cells.Item(31, 6).Value = 1203  # F31: 1201 -> 1203
$ws4.Cells.Item(33, 6).Value = 902  # F33: 898 -> 902
$ws4.Cells.Item(34, 6).Value = 2630  # F34: 2624 -> 2630
$ws4.Cells.Item(35, 6).Value = 1120  # F35: 1111 -> 1120
$ws4.Cells.Item(36, 6).Value = 592  # F36: 590 -> 592
$ws4.Cells.Item(37, 6).Value = 134  # F37: 133 -> 134
$ws4.Cells.Item(40, 6).Value = 16  # F40: 14 -> 16
$ws4.Cells.Item(41, 6).Value = 867  # F41: 863 -> 867
$ws4.Cells.Item(43, 6).Value = 1206  # F43: 1202 -> 1206
